# Actualiza el periodo de mora de "2508" a "2509" para todos los trabajadores
# y centra horizontalmente la columna "Periodo Mora" (E16:E21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Periodo Mora" values (column E, rows 16-21) from 2508 to 2509
$ws.Range("E16:E21").Value = "2509"

# Center-align the "Periodo Mora" column values
$ws.Range("E16:E21").HorizontalAlignment = -4108
